$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.781.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.02%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.825.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.96%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'615.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'177.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.43%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.823.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.99%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.29%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.32%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.81%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'40.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.60%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.452.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.92%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.814.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.66%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'69.879.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.10%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'7.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.61%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -3.65%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'16.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.38%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'511.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D23").Value = "'0.740"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.38%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.82%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'86.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.72%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0000145"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.27%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.22%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.28%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.20%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.50%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +3.45%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +2.48%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'31.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.58%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.06%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.72%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +7.68%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'485.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +14.05%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.45%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.92%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.96%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'49.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.26%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'44.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.56%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.34%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.959.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.50%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0365"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.11%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'27.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.59%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'139.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.96%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.63%  "
$ws.Range("E51").Style = "Normal"
